# This workbook contains a weekly price log for "Espinaca" (spinach) at the
# "Vega Modelo de Temuco" market. A new weekly record was added; it is
# inserted as a new row 57 (pushing the existing rows 57-127 down to 58-128),
# matching how the rest of the table is ordered.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 57 - this shifts rows 57..127 down to
# 58..128 (and grows the used range to A1:R128), which already reproduces
# all of the "shifted" values seen in the diff for rows 58-128.
$ws.Rows.Item(57).Insert()

# Populate the newly inserted row 57 with the new weekly record.
$ws.Cells.Item(57, 1).Value  = 10
$ws.Cells.Item(57, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(57, 3).Value  = "La Araucanía"
$ws.Cells.Item(57, 4).Value  = 44638
$ws.Cells.Item(57, 5).Value  = 9
$ws.Cells.Item(57, 6).Value  = 100112012
$ws.Cells.Item(57, 7).Value  = "Espinaca"
$ws.Cells.Item(57, 8).Value  = "Sin especificar"
$ws.Cells.Item(57, 9).Value  = "Primera"
$ws.Cells.Item(57, 10).Value = 20
$ws.Cells.Item(57, 11).Value = 10000
$ws.Cells.Item(57, 12).Value = 10000
$ws.Cells.Item(57, 13).Value = 10000
$ws.Cells.Item(57, 14).Value = "$/docena de atados"
$ws.Cells.Item(57, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(57, 16).Value = 3333
$ws.Cells.Item(57, 17).Value = 3
$ws.Cells.Item(57, 18).Value = "Hortaliza"
